$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.4529906666666667
$ws.Range("H2").Value = 1.358972
$ws.Range("I2").Value = 0.1117470803109675
$ws.Range("J2").Value = 0.1238493010025621
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.638706
$ws.Range("N2").Value = 1.916118
$ws.Range("O2").Value = 0.5586654432763536
$ws.Range("P2").Value = 0.5586654432763536
$ws.Range("Q2").Value = 0.289327856744
$ws.Range("R2").Value = 2.603950710696
$ws.Range("S2").Value = 0.06242923215676496
$ws.Range("T2").Value = 0.06919032464406288

$ws.Range("G3").Value = 0.4529906666666667
$ws.Range("H3").Value = 1.358972
$ws.Range("I3").Value = 0.1117470803109675
$ws.Range("J3").Value = 0.1238493010025621
$ws.Range("M3").Value = 0.5045649999999999
$ws.Range("N3").Value = 1.513695
$ws.Range("O3").Value = 0.4413345567236464
$ws.Range("P3").Value = 0.4413345567236464
$ws.Range("Q3").Value = 0.2285632357266666
$ws.Range("R3").Value = 2.05706912154
$ws.Range("S3").Value = 0.04931784815420258
$ws.Range("T3").Value = 0.05465897635849919

$ws.Range("H4").Value = 5.124729
$ws.Range("I4").Value = 0.4214019885140712
$ws.Range("J4").Value = 0.4670398687225041
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.638706
$ws.Range("N4").Value = 1.916118
$ws.Range("O4").Value = 0.5586654432763536
$ws.Range("P4").Value = 0.5586654432763536
$ws.Range("Q4").Value = 1.091065053558
$ws.Range("R4").Value = 9.819585482022001
$ws.Range("S4").Value = 0.2354227287107504
$ws.Range("T4").Value = 0.2609190352875878

$ws.Range("H5").Value = 5.124729
$ws.Range("I5").Value = 0.4214019885140712
$ws.Range("J5").Value = 0.4670398687225041
$ws.Range("M5").Value = 0.5045649999999999
$ws.Range("N5").Value = 1.513695
$ws.Range("O5").Value = 0.4413345567236464
$ws.Range("P5").Value = 0.4413345567236464
$ws.Range("Q5").Value = 0.861919629295
$ws.Range("R5").Value = 7.757276663654999
$ws.Range("S5").Value = 0.1859792598033208
$ws.Range("T5").Value = 0.2061208334349164

$ws.Range("G6").Value = 0.1307166666666667
$ws.Range("H6").Value = 0.39215
$ws.Range("I6").Value = 0.03224615190301634
$ws.Range("J6").Value = 0.03573841358626573
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.638706
$ws.Range("N6").Value = 1.916118
$ws.Range("O6").Value = 0.5586654432763536
$ws.Range("P6").Value = 0.5586654432763536
$ws.Range("Q6").Value = 0.08348951930000001
$ws.Range("R6").Value = 0.7514056737
$ws.Range("S6").Value = 0.01801481074685525
$ws.Range("T6").Value = 0.01996581666816481

$ws.Range("G7").Value = 0.1307166666666667
$ws.Range("H7").Value = 0.39215
$ws.Range("I7").Value = 0.03224615190301634
$ws.Range("J7").Value = 0.03573841358626573
$ws.Range("M7").Value = 0.5045649999999999
$ws.Range("N7").Value = 1.513695
$ws.Range("O7").Value = 0.4413345567236464
$ws.Range("P7").Value = 0.4413345567236464
$ws.Range("Q7").Value = 0.06595505491666666
$ws.Range("R7").Value = 0.5935954942499999
$ws.Range("S7").Value = 0.01423134115616108
$ws.Range("T7").Value = 0.01577259691810093

$ws.Range("G8").Value = 1.188354
$ws.Range("H8").Value = 2.376708
$ws.Range("I8").Value = 0.2931519336878011
$ws.Range("J8").Value = 0.2166002128720807
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.638706
$ws.Range("N8").Value = 1.916118
$ws.Range("O8").Value = 0.5586654432763536
$ws.Range("P8").Value = 0.5586654432763536
$ws.Range("Q8").Value = 0.759008829924
$ws.Range("R8").Value = 4.554052979543999
$ws.Range("S8").Value = 0.1637738549810156
$ws.Range("T8").Value = 0.1210070539379335

$ws.Range("G9").Value = 1.188354
$ws.Range("H9").Value = 2.376708
$ws.Range("I9").Value = 0.2931519336878011
$ws.Range("J9").Value = 0.2166002128720807
$ws.Range("M9").Value = 0.5045649999999999
$ws.Range("N9").Value = 1.513695
$ws.Range("O9").Value = 0.4413345567236464
$ws.Range("P9").Value = 0.4413345567236464
$ws.Range("Q9").Value = 0.5996018360099998
$ws.Range("R9").Value = 3.597611016059999
$ws.Range("S9").Value = 0.1293780787067855
$ws.Range("T9").Value = 0.09559315893414719

$ws.Range("G10").Value = 0.5734093333333333
$ws.Range("H10").Value = 1.720228
$ws.Range("I10").Value = 0.1414528455841438
$ws.Range("J10").Value = 0.1567722038165874
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.638706
$ws.Range("N10").Value = 1.916118
$ws.Range("O10").Value = 0.5586654432763536
$ws.Range("P10").Value = 0.5586654432763536
$ws.Range("Q10").Value = 0.366239981656
$ws.Range("R10").Value = 3.296159834904
$ws.Range("S10").Value = 0.07902481668096728
$ws.Range("T10").Value = 0.08758321273860463

$ws.Range("G11").Value = 0.5734093333333333
$ws.Range("H11").Value = 1.720228
$ws.Range("I11").Value = 0.1414528455841438
$ws.Range("J11").Value = 0.1567722038165874
$ws.Range("M11").Value = 0.5045649999999999
$ws.Range("N11").Value = 1.513695
$ws.Range("O11").Value = 0.4413345567236464
$ws.Range("P11").Value = 0.4413345567236464
$ws.Range("Q11").Value = 0.2893222802733333
$ws.Range("R11").Value = 2.60390052246
$ws.Range("S11").Value = 0.0624280289031765
$ws.Range("T11").Value = 0.06918899107798272
